# Update cryptos list: apply new Price / Volume(1h) figures and the
# Dai / InternetComputer(DFINITY) row swap (rows 30-31).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text is not a "clean" number (so Excel's COM layer will
# naturally keep them as text, matching the original inlineStr cells).
$plainUpdates = @{
    'D2' = '93.045.64'
    'E2' = '  -5.31%  '
    'D3' = '3.372.29'
    'E3' = '  -1.25%  '
    'E4' = '  +0.04%  '
    'E5' = '  -8.92%  '
    'E6' = '  -6.99%  '
    'E7' = '  -8.60%  '
    'E8' = '  -10.40%  '
    'E9' = '  +0.13%  '
    'E10' = '  -12.82%  '
    'D11' = '3.374.30'
    'E11' = '  -1.11%  '
    'E12' = '  -7.24%  '
    'E13' = '  -13.35%  '
    'E14' = '  -2.60%  '
    'D15' = '92.940.07'
    'E15' = '  -5.23%  '
    'D16' = '3.991.45'
    'E16' = '  -1.94%  '
    'E17' = '  -6.52%  '
    'E18' = '  -12.53%  '
    'D19' = '3.367.27'
    'E19' = '  -1.54%  '
    'E20' = '  -8.82%  '
    'E21' = '  -5.26%  '
    'E22' = '  -5.43%  '
    'E23' = '  -16.09%  '
    'E24' = '  -9.45%  '
    'E25' = '  -8.67%  '
    'E26' = '  -8.69%  '
    'E27' = '  -8.44%  '
    'D28' = '3.539.85'
    'E28' = '  -1.45%  '
    'E29' = '  -8.83%  '
    'B30' = 'Dai'
    'C30' = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
    'E30' = '  +0.08%  '
    'B31' = 'InternetComputer(DFINITY)'
    'C31' = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
    'E31' = '  -8.61%  '
    'E32' = '  -12.59%  '
    'E33' = '  -10.07%  '
    'E34' = '  -1.31%  '
    'E35' = '  -10.41%  '
    'E36' = '  -2.79%  '
    'E37' = '  -7.67%  '
    'E38' = '  +0.35%  '
    'E39' = '  -7.40%  '
    'E40' = '  -0.08%  '
    'E41' = '  -7.36%  '
    'E42' = '  -6.02%  '
    'E43' = '  +0.42%  '
    'E44' = '  -1.67%  '
    'E45' = '  -2.57%  '
    'E46' = '  -7.10%  '
    'E47' = '  -4.25%  '
    'E48' = '  -5.03%  '
    'E49' = '  -5.09%  '
    'E50' = '  -9.43%  '
    'E51' = '  -3.91%  '
}

foreach ($addr in $plainUpdates.Keys) {
    $ws.Range($addr).Value = $plainUpdates[$addr]
}

# Cells whose new text DOES look like a clean number (e.g. "1.00", "0.130").
# Assigning those directly would make Excel coerce the cell to a Number
# (and drop the formatting, e.g. trailing zeros), so force them to stay
# text the same way typing `'1.00` in Excel does, then restore the
# default "Normal" style so no stray number-format style is left behind.
$numericUpdates = @{
    'D5' = '232.14'
    'D6' = '625.85'
    'D8' = '0.387'
    'D10' = '0.928'
    'D13' = '40.17'
    'D17' = '0.0000242'
    'D18' = '7.97'
    'D20' = '16.83'
    'D21' = '10.96'
    'D22' = '485.74'
    'D23' = '0.454'
    'D26' = '6.25'
    'D27' = '89.66'
    'D30' = '1.00'
    'D31' = '11.25'
    'D33' = '0.130'
    'D34' = '0.984'
    'D35' = '0.169'
    'D36' = '28.49'
    'D37' = '0.529'
    'D38' = '533.74'
    'D39' = '7.43'
    'D43' = '0.873'
    'D44' = '24.01'
    'D45' = '3.59'
    'D46' = '1.64'
    'D47' = '5.45'
    'D49' = '53.01'
    'D50' = '0.0389'
    'D51' = '3.13'
}

foreach ($addr in $numericUpdates.Keys) {
    $ws.Range($addr).Value = "'" + $numericUpdates[$addr]
    $ws.Range($addr).Style = "Normal"
}
